$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (MAE)
$ws.Range("B2").Value = 2.235
$ws.Range("C2").Value = 2.09
$ws.Range("D2").Value = 1.773
$ws.Range("E2").Value = 1.321
$ws.Range("F2").Value = 2.371

# Row 3 (MSE)
$ws.Range("B3").Value = 6.359
$ws.Range("C3").Value = 7.11
$ws.Range("D3").Value = 5.37
$ws.Range("E3").Value = 2.588
$ws.Range("F3").Value = 8.487

# Row 4 (mean Y-Test)
$ws.Range("B4").Value = 18.203
$ws.Range("C4").Value = 15.45
$ws.Range("D4").Value = 18.059
$ws.Range("E4").Value = 12.974
$ws.Range("F4").Value = 30.217

# Row 5 (mean Y-predicted)
$ws.Range("B5").Value = 18.5
$ws.Range("C5").Value = 15.768
$ws.Range("D5").Value = 18.055
$ws.Range("E5").Value = 13.627
$ws.Range("F5").Value = 31.057

# Row 6 (R2)
$ws.Range("B6").Value = 0.47
$ws.Range("C6").Value = 0.556
$ws.Range("D6").Value = 0.772
$ws.Range("E6").Value = 0.228
$ws.Range("F6").Value = 0.783
